# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have been generated: the Status column flips
# from "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns are populated on both locale sheets, hyperlinks are
# added for the new "Latest Target File" links, and a few columns are
# widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e1ca09ac2bd913f8b030dac40b2055571fed10/e2e/a.md"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared across the Overview summary columns and both locale sheets)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus

$zhcn.Range("C2").Value2 = $newStatus
$zhcn.Range("C3").Value2 = $newStatus

$dede.Range("C2").Value2 = $newStatus
$dede.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------
# 2. Populate "Latest Handback File" (J) / "Latest Handback DateTime" (K)
#    for both locale sheets.
# ---------------------------------------------------------------------
$zhcn.Range("J2").Value2 = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value2 = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value2 = "2016-08-31 00:39:49"
$zhcn.Range("K3").Value2 = "2016-08-31 00:39:49"

$dede.Range("J2").Value2 = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value2 = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value2 = "2016-08-31 00:39:56"
$dede.Range("K3").Value2 = "2016-08-31 00:39:56"

# ---------------------------------------------------------------------
# 3. Populate "Latest Target File" (I) with a hyperlink to a.md on both
#    rows of both locale sheets. The existing hyperlinks (A2/A3) are
#    re-created alongside so that the full set ends up in document order
#    (A2, I2, A3, I3) exactly like Excel lays out a freshly-saved sheet.
# ---------------------------------------------------------------------
function Set-TargetFileHyperlinks($ws, $bDisplay) {
    $addrA2 = ""
    $dispA2 = ""
    $addrA3 = ""
    $dispA3 = ""
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq 2) {
            $addrA2 = $hl.Address
            $dispA2 = $hl.TextToDisplay
        } elseif ($hl.Range.Row -eq 3) {
            $addrA3 = $hl.Address
            $dispA3 = $hl.TextToDisplay
        }
    }

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $addrA2, [Type]::Missing, [Type]::Missing, $dispA2)
    $ws.Hyperlinks.Add($ws.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $addrA3, [Type]::Missing, [Type]::Missing, $dispA3)
    $ws.Hyperlinks.Add($ws.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
}

Set-TargetFileHyperlinks $zhcn "b.md"
Set-TargetFileHyperlinks $dede "b.md"

# ---------------------------------------------------------------------
# 4. Widen columns to fit the new, longer content.
# ---------------------------------------------------------------------
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("J1").ColumnWidth = 40

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("J1").ColumnWidth = 40
